$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lugar")

# Insert three new columns before the existing F/G (img-stand / img-comedero),
# which pushes those two columns to I/J.
$ws.Range("F1:H1").EntireColumn.Insert()

# New header row values (F1, G1, H1). Shared strings must be created in the
# order Descripcion, Cap. Comedor (kg), Cap. Stand (Personas), so populate
# H1 before G1.
$ws.Cells.Item(1, 6).Value = "Descripcion"
$ws.Cells.Item(1, 8).Value = "Cap. Comedor (kg)"
$ws.Cells.Item(1, 7).Value = "Cap. Stand (Personas)"

# Match the new column widths from the target layout as closely as this
# quantized COM property allows.
$ws.Range("F1:H1").EntireColumn.ColumnWidth = 19.6
$ws.Range("I1:I1").EntireColumn.ColumnWidth = 17.5

# Data rows: Descripcion / Cap. Stand (Personas) / Cap. Comedor (kg).
$descripciones = @("Blind", "Blind", "Blind", "Rancho", "Blind", "Blind", "Blind")
for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $descripciones[$i]
    $ws.Cells.Item($row, 7).Value = 2
    $ws.Cells.Item($row, 8).Value = 95
}

# Make "Lugar" the active sheet/tab (it was "Hitorial" before) and select
# the new active cell on it; "Hitorial" keeps its existing E6 selection.
[void]$ws.Activate()
[void]$ws.Range("H14").Select()
